# Apply the "cryptos list" update described by the commit diff.
# - Rows 2-37: price (D) / volume-change (E) values are refreshed in place.
# - A new "BinanceUSD" row is inserted at row 38, pushing the former rows
#   38-50 down to 39-51 and dropping the old trailing row 51 (TrustWalletToken).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the literal text into the cell even when it looks numeric
    # (e.g. "1.00", "305.85"), then drop back to the default style so we
    # do not leave a stray text-format override behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---- Rows 2-37: refresh Price (D) / Volume(1h) (E) cells ----
Set-TextValue $ws.Range('D2') '42.045.15'
$ws.Range('E2').Value = '  +0.32%  '
Set-TextValue $ws.Range('D3') '2.269.62'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '305.85'
$ws.Range('E5').Value = '  +1.49%  '
Set-TextValue $ws.Range('D6') '93.14'
$ws.Range('E6').Value = '  +1.40%  '
Set-TextValue $ws.Range('D7') '0.530'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.63%  '
Set-TextValue $ws.Range('D10') '32.84'
$ws.Range('E10').Value = '  +1.97%  '
Set-TextValue $ws.Range('D11') '0.0801'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('E13').Value = '  +0.77%  '
Set-TextValue $ws.Range('D14') '2.622.54'
$ws.Range('E14').Value = '  +0.63%  '
Set-TextValue $ws.Range('D15') '14.35'
$ws.Range('E15').Value = '  +1.93%  '
Set-TextValue $ws.Range('D16') '2.285.98'
$ws.Range('E16').Value = '  +1.47%  '
Set-TextValue $ws.Range('D17') '0.785'
$ws.Range('E17').Value = '  +4.02%  '
Set-TextValue $ws.Range('D18') '41.917.14'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').Value = '  +5.78%  '
$ws.Range('E20').Value = '  +2.06%  '
$ws.Range('E21').Value = '  +1.34%  '
Set-TextValue $ws.Range('D22') '68.19'
$ws.Range('E22').Value = '  +1.90%  '
Set-TextValue $ws.Range('D23') '244.00'
$ws.Range('E23').Value = '  +1.46%  '
Set-TextValue $ws.Range('D24') '2.60'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('E26').Value = '  +0.07%  '
Set-TextValue $ws.Range('D27') '23.98'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('E29').Value = '  -9.06%  '
Set-TextValue $ws.Range('D30') '34.89'
$ws.Range('E30').Value = '  +3.61%  '
Set-TextValue $ws.Range('D31') '159.85'
$ws.Range('E31').Value = '  +1.04%  '
Set-TextValue $ws.Range('D32') '5.34'
$ws.Range('E32').Value = '  +4.05%  '
$ws.Range('E33').Value = '  +0.01%  '
Set-TextValue $ws.Range('D34') '0.0744'
$ws.Range('E34').Value = '  +0.56%  '
Set-TextValue $ws.Range('D35') '3.03'
$ws.Range('E35').Value = '  +0.10%  '
Set-TextValue $ws.Range('D36') '17.14'
$ws.Range('E36').Value = '  +4.54%  '
$ws.Range('E37').Value = '  -1.00%  '

# ---- Insert the new BinanceUSD row at 38 ----
# This shifts the old rows 38-50 down to 39-51 and pushes what was
# row 51 (TrustWalletToken) off the bottom of the A1:E51 range.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(52).Delete()

# Column A keeps the bold / centered / thin-border header style used by
# every other row; the freshly inserted row does not inherit it automatically.
$aNew = $ws.Range("A38")
$aNew.Font.Bold = $true
$aNew.HorizontalAlignment = -4108
$aNew.VerticalAlignment = -4160
$aNew.Borders.LineStyle = 1

# ---- Rows 38-51: write the post-shift data ----

# row 38: BinanceUSD
$ws.Range('A38').Value = 36
$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D38') '4.42'
$ws.Range('E38').Value = '  +341.43%  '

# row 39: Kaspa
$ws.Range('A39').Value = 37
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D39') '0.105'
$ws.Range('E39').Value = '  +1.40%  '

# row 40: Stellar
$ws.Range('A40').Value = 38
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D40') '0.117'
$ws.Range('E40').Value = '  +1.20%  '

# row 41: ARBITRUM
$ws.Range('A41').Value = 39
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D41') '1.79'
$ws.Range('E41').Value = '  +0.73%  '

# row 42: RenderToken
$ws.Range('A42').Value = 40
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D42') '3.99'
$ws.Range('E42').Value = '  +1.94%  '

# row 43: EnergySwap
$ws.Range('A43').Value = 41
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D43') '19.83'
$ws.Range('E43').Value = '  +0.33%  '

# row 44: Maker
$ws.Range('A44').Value = 42
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D44') '2.016.89'
$ws.Range('E44').Value = '  -1.54%  '

# row 45: ApeXProtocol
$ws.Range('A45').Value = 43
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range('D45') '2.24'
$ws.Range('E45').Value = '  +9.02%  '

# row 46: VeChain
$ws.Range('A46').Value = 44
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D46') '0.0283'
$ws.Range('E46').Value = '  +1.79%  '

# row 47: FraxShare
$ws.Range('A47').Value = 45
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D47') '10.28'
$ws.Range('E47').Value = '  +2.09%  '

# row 48: NEARProtocol
$ws.Range('A48').Value = 46
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D48') '2.91'
$ws.Range('E48').Value = '  +3.20%  '

# row 49: MultiversX
$ws.Range('A49').Value = 47
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range('D49') '53.31'
$ws.Range('E49').Value = '  +3.31%  '

# row 50: Stacks
$ws.Range('A50').Value = 48
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D50') '1.52'
$ws.Range('E50').Value = '  +0.69%  '

# row 51: BitcoinSV
$ws.Range('A51').Value = 49
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range('D51') '72.46'
$ws.Range('E51').Value = '  +3.08%  '
